# Update the "WEIGHT ESTIMATION METHODS COMPARISON" tables across several
# sheets of the JPAD Weights workbook. The underlying shared-strings table
# for the method-name labels (NICOLAI_1984, SADRAEY, JENKINSON, ROSKAM,
# KROO, RAYMER, TORENBEEK_2013, TORENBEEK_1976, TORENBEEK_1982, HOWE,
# NICOLAI_2013, ...) was reordered by the report generator, and the rows
# that reference those labels were refreshed so each row keeps showing the
# correct method name together with its correct Estimated Mass / Percent
# Error figures.
#
# Rather than poke at shared-string indices directly (an implementation
# detail Excel manages on its own), we just (re)write the final displayed
# values for every affected cell on each sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("A8").Value = "NICOLAI_1984"
$ws.Range("C8").Value = 2968.0
$ws.Range("D8").Value = 9.169624360409161
$ws.Range("A9").Value = "SADRAEY"
$ws.Range("C9").Value = 2491.0
$ws.Range("D9").Value = -8.375493840370883
$ws.Range("A10").Value = "JENKINSON"
$ws.Range("C10").Value = 4506.0
$ws.Range("D10").Value = 65.74067633692846
$ws.Range("A11").Value = "ROSKAM"
$ws.Range("C11").Value = 3917.0
$ws.Range("D11").Value = 44.07594966971788
$ws.Range("A12").Value = "KROO"
$ws.Range("C12").Value = 2585.0
$ws.Range("D12").Value = -4.917965306045256
$ws.Range("A13").Value = "RAYMER"
$ws.Range("C13").Value = 3149.0
$ws.Range("D13").Value = 15.827205899908506
$ws.Range("A14").Value = "TORENBEEK_2013"
$ws.Range("C14").Value = 3698.0
$ws.Range("D14").Value = 36.020643829108174
$ws.Range("A15").Value = "TORENBEEK_1976"
$ws.Range("C15").Value = 3818.0
$ws.Range("D15").Value = 40.43451004314089

$ws = $wb.Worksheets.Item("WING")
$ws.Range("A9").Value = "KROO"
$ws.Range("C9").Value = 2539.0
$ws.Range("D9").Value = -11.89617674977151
$ws.Range("A10").Value = "RAYMER"
$ws.Range("C10").Value = 2760.0
$ws.Range("D10").Value = -4.227431204950519
$ws.Range("A11").Value = "TORENBEEK_2013"

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("A8").Value = "SADRAEY"
$ws.Range("C8").Value = 273.0
$ws.Range("D8").Value = -12.68221185283106
$ws.Range("A9").Value = "JENKINSON"
$ws.Range("C9").Value = 293.0
$ws.Range("D9").Value = -6.285304296261907
$ws.Range("A10").Value = "HOWE"
$ws.Range("C10").Value = 207.0
$ws.Range("D10").Value = -33.792006789509266
$ws.Range("A11").Value = "NICOLAI_2013"
$ws.Range("C11").Value = 124.0
$ws.Range("D11").Value = -60.339173149271254
$ws.Range("A12").Value = "KROO"
$ws.Range("C12").Value = 303.0
$ws.Range("D12").Value = -3.0868505179773305
$ws.Range("A13").Value = "ROSKAM"
$ws.Range("C13").Value = 216.0
$ws.Range("D13").Value = -30.91339838905315
$ws.Range("A14").Value = "RAYMER"
$ws.Range("C14").Value = 144.0
$ws.Range("D14").Value = -53.9422655927021
$ws.Range("A15").Value = "TORENBEEK_1976"
$ws.Range("C15").Value = 236.0
$ws.Range("D15").Value = -24.516490832483992

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("A8").Value = "SADRAEY"
$ws.Range("C8").Value = 413.0
$ws.Range("D8").Value = 32.09614104315301
$ws.Range("A9").Value = "JENKINSON"
$ws.Range("C9").Value = 330.0
$ws.Range("D9").Value = 5.548974683391027
$ws.Range("A10").Value = "HOWE"
$ws.Range("C10").Value = 445.0
$ws.Range("D10").Value = 42.33119313366365
$ws.Range("A11").Value = "ROSKAM"
$ws.Range("C11").Value = 239.0
$ws.Range("D11").Value = -23.55695469899862
$ws.Range("A12").Value = "KROO"
$ws.Range("C12").Value = 256.0
$ws.Range("D12").Value = -18.11958327591484
$ws.Range("A13").Value = "RAYMER"
$ws.Range("C13").Value = 89.0
$ws.Range("D13").Value = -71.53376137326727
$ws.Range("A14").Value = "TORENBEEK_1976"
$ws.Range("C14").Value = 338.0
$ws.Range("D14").Value = 8.107737706018687

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("A11").Value = "TORENBEEK_1976"
$ws.Range("A17").Value = "TORENBEEK_1976"

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("A12").Value = "TORENBEEK_1976"
$ws.Range("A18").Value = "TORENBEEK_1976"

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("A9").Value = "TORENBEEK_1976"
$ws.Range("A11").Value = "TORENBEEK_1976"
$ws.Range("A13").Value = "TORENBEEK_1976"
